# Refresh cryptocurrency price / volume snapshot (column D) and two
# "Worstin24h" label toggles (column E) to match the latest scrape,
# as produced by the scheduled GitHub Actions symbol-list update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the cell address together with its new literal text.
# All of these cells are plain text cells in the source sheet (prices
# such as 0.0005824 must keep their exact textual formatting, i.e. not
# be reinterpreted as floating point numbers), so for every cell we
# force text formatting ("@") before assigning the value and then drop
# back to the workbook's default (Normal) style so no stray per-cell
# formatting is introduced.
$updates = @(
    @{Cell="D2"; Value="250.68"},
    @{Cell="D3"; Value="21.78"},
    @{Cell="D4"; Value="5.561"},
    @{Cell="D5"; Value="0.05689"},
    @{Cell="D6"; Value="6.450"},
    @{Cell="D7"; Value="0.8078"},
    @{Cell="D8"; Value="1.042"},
    @{Cell="D9"; Value="0.1432"},
    @{Cell="D10"; Value="0.07277"},
    @{Cell="D11"; Value="0.03136"},
    @{Cell="D12"; Value="0.02923"},
    @{Cell="D13"; Value="0.09278"},
    @{Cell="D14"; Value="0.001651"},
    @{Cell="D15"; Value="3.228"},
    @{Cell="D16"; Value="0.04754"},
    @{Cell="D17"; Value="0.0005824"},
    @{Cell="E17"; Value="16OneONEWorstin24h"},
    @{Cell="D18"; Value="0.006459"},
    @{Cell="D19"; Value="0.005065"},
    @{Cell="D20"; Value="0.001055"},
    @{Cell="D21"; Value="0.0001502"},
    @{Cell="D22"; Value="3.986"},
    @{Cell="D23"; Value="3.377"},
    @{Cell="D24"; Value="2.114"},
    @{Cell="D25"; Value="0.3321"},
    @{Cell="D27"; Value="0.0003102"},
    @{Cell="D40"; Value="0.04126"},
    @{Cell="D41"; Value="0.006908"},
    @{Cell="D42"; Value="0.1049"},
    @{Cell="D43"; Value="0.003203"},
    @{Cell="D44"; Value="0.009560"},
    @{Cell="D45"; Value="0.00005643"},
    @{Cell="D46"; Value="0.00000000750"},
    @{Cell="D47"; Value="0.7856"},
    @{Cell="D48"; Value="0.01690"},
    @{Cell="E48"; Value="47BOLOBOLO"},
    @{Cell="D49"; Value="0.00002101"},
    @{Cell="D50"; Value="0.01010"}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
